$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31 (cohort 2022, period_index 3): num_customers 41 -> 42, retention_rate recalculated
$ws.Range("C31").Value = 42
$ws.Range("E31").Value = 42 / 2312

# Row 34 (cohort 2023, period_index 2): num_customers 71 -> 72, retention_rate recalculated
$ws.Range("C34").Value = 72
$ws.Range("E34").Value = 72 / 2256

# Row 36 (cohort 2024, period_index 1): num_customers 114 -> 115, retention_rate recalculated
$ws.Range("C36").Value = 115
$ws.Range("E36").Value = 115 / 1930

# Row 37 (cohort 2025, period_index 0): num_customers 713 -> 720, cohort_size 713 -> 720
$ws.Range("C37").Value = 720
$ws.Range("D37").Value = 720
